# Datadrivern [property files, Csvfiles]
#
# The "Sign in" button locator on the Outlook-validation sheet (Sheet5) is
# updated to a more robust XPath expression, and that sheet is brought to
# the foreground (becomes the active/selected tab), with the cursor left on
# cell B8 - mirroring the state Excel saves after a user clicks into that
# sheet and edits the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet5")

# Bring Sheet5 to the front (updates workbook activeTab / tabSelected).
$ws.Activate()

# Update the Sign-in button xpath locator.
$ws.Range("B7").Value = "(//a[contains(.,'Sign in')])[1]"

# Leave the selection on B8, as recorded in the saved view state.
$ws.Range("B8").Select()
